# Auto-generated edit script: updates cached market-price / profit values
# across the 8 Leve profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# reflecting a refreshed data pull. One cell (GSM!N94) is cleared entirely
# because its source row no longer yields an HQ profit value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 188.08824
$ws.Range("I9").Value = 322.7143
$ws.Range("K9").Value = 322.7143
$ws.Range("M9").Value = -153.7143
$ws.Range("H40").Value = 4224.615
$ws.Range("J40").Value = 5279.2856
$ws.Range("L40").Value = 5279.2856
$ws.Range("N40").Value = -5629.2856
$ws.Range("H62").Value = 3511.0833
$ws.Range("I62").Value = 3154.2
$ws.Range("J62").Value = 4105.8887
$ws.Range("K62").Value = 3154.2
$ws.Range("L62").Value = 4105.8887
$ws.Range("M62").Value = -2530.2
$ws.Range("N62").Value = -5353.8887
$ws.Range("H65").Value = 3511.0833
$ws.Range("I65").Value = 3154.2
$ws.Range("J65").Value = 4105.8887
$ws.Range("K65").Value = 15771
$ws.Range("L65").Value = 20529.4435
$ws.Range("M65").Value = -12651
$ws.Range("N65").Value = -26769.4435
$ws.Range("H76").Value = 5743.1816
$ws.Range("I76").Value = 5758.4443
$ws.Range("K76").Value = 5758.4443
$ws.Range("M76").Value = -5443.4443
$ws.Range("H79").Value = 5743.1816
$ws.Range("I79").Value = 5758.4443
$ws.Range("K79").Value = 5758.4443
$ws.Range("M79").Value = -4666.4443
$ws.Range("H93").Value = 49999.5
$ws.Range("J93").Value = 49999.5
$ws.Range("L93").Value = 49999.5
$ws.Range("N93").Value = -54991.5
$ws.Range("H96").Value = 762.2941
$ws.Range("I96").Value = 467.4
$ws.Range("K96").Value = 1402.2
$ws.Range("M96").Value = -29.19999999999982
$ws.Range("H98").Value = 1290.3793
$ws.Range("I98").Value = 706.7619
$ws.Range("J98").Value = 2822.375
$ws.Range("K98").Value = 706.7619
$ws.Range("L98").Value = 2822.375
$ws.Range("M98").Value = 791.2381
$ws.Range("N98").Value = -5818.375
$ws.Range("H106").Value = 20330
$ws.Range("I106").Value = 3330
$ws.Range("K106").Value = 3330
$ws.Range("M106").Value = -2699
$ws.Range("H122").Value = 1290.3793
$ws.Range("I122").Value = 706.7619
$ws.Range("J122").Value = 2822.375
$ws.Range("K122").Value = 2120.2857
$ws.Range("L122").Value = 8467.125
$ws.Range("M122").Value = 329.7143000000001
$ws.Range("N122").Value = -13367.125
$ws.Range("H131").Value = 743.4545000000001
$ws.Range("I131").Value = 656
$ws.Range("J131").Value = 1137
$ws.Range("K131").Value = 1968
$ws.Range("L131").Value = 3411
$ws.Range("M131").Value = 3072
$ws.Range("N131").Value = -13491
$ws.Range("H132").Value = 3011.3684
$ws.Range("I132").Value = 2713.625
$ws.Range("K132").Value = 8140.875
$ws.Range("M132").Value = -5610.875
$ws.Range("H137").Value = 1234.2307
$ws.Range("I137").Value = 1060.25
$ws.Range("J137").Value = 1814.1666
$ws.Range("K137").Value = 3180.75
$ws.Range("L137").Value = 5442.4998
$ws.Range("M137").Value = -630.75
$ws.Range("N137").Value = -10542.4998
$ws.Range("H138").Value = 2603.2903
$ws.Range("I138").Value = 3674.889
$ws.Range("J138").Value = 2164.9092
$ws.Range("K138").Value = 11024.667
$ws.Range("L138").Value = 6494.7276
$ws.Range("M138").Value = -5884.667000000001
$ws.Range("N138").Value = -16774.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 68
$ws.Range("I5").Value = 57.666668
$ws.Range("J5").Value = 99
$ws.Range("K5").Value = 57.666668
$ws.Range("L5").Value = 99
$ws.Range("M5").Value = 54.333332
$ws.Range("N5").Value = -323
$ws.Range("H32").Value = 28305.13
$ws.Range("I32").Value = 6280.073
$ws.Range("K32").Value = 6280.073
$ws.Range("M32").Value = -5993.073
$ws.Range("H45").Value = 5008.657
$ws.Range("I45").Value = 5860.913
$ws.Range("K45").Value = 5860.913
$ws.Range("M45").Value = -5483.913
$ws.Range("H52").Value = 24917.4
$ws.Range("J52").Value = 24917.4
$ws.Range("L52").Value = 24917.4
$ws.Range("N52").Value = -25553.4
$ws.Range("H61").Value = 2370.9583
$ws.Range("I61").Value = 2071.45
$ws.Range("K61").Value = 2071.45
$ws.Range("M61").Value = -1859.45
$ws.Range("H74").Value = 2076.8667
$ws.Range("I74").Value = 1836.88
$ws.Range("J74").Value = 3276.8
$ws.Range("K74").Value = 1836.88
$ws.Range("L74").Value = 3276.8
$ws.Range("M74").Value = -962.8800000000001
$ws.Range("N74").Value = -5024.8
$ws.Range("H77").Value = 2076.8667
$ws.Range("I77").Value = 1836.88
$ws.Range("J77").Value = 3276.8
$ws.Range("K77").Value = 9184.400000000001
$ws.Range("L77").Value = 16384
$ws.Range("M77").Value = -4816.400000000001
$ws.Range("N77").Value = -25120
$ws.Range("H102").Value = 3323.9412
$ws.Range("I102").Value = 1398.8462
$ws.Range("K102").Value = 1398.8462
$ws.Range("M102").Value = 223.1538
$ws.Range("H122").Value = 2686.5908
$ws.Range("I122").Value = 2689.1765
$ws.Range("J122").Value = 2677.8
$ws.Range("K122").Value = 8067.529500000001
$ws.Range("L122").Value = 8033.400000000001
$ws.Range("M122").Value = -5617.529500000001
$ws.Range("N122").Value = -12933.4
$ws.Range("H136").Value = 2370.9583
$ws.Range("I136").Value = 2071.45
$ws.Range("K136").Value = 6214.349999999999
$ws.Range("M136").Value = -3664.349999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 68
$ws.Range("I4").Value = 57.666668
$ws.Range("J4").Value = 99
$ws.Range("K4").Value = 57.666668
$ws.Range("L4").Value = 99
$ws.Range("M4").Value = 57.333332
$ws.Range("N4").Value = -329
$ws.Range("H20").Value = 10218.84
$ws.Range("I20").Value = 8078.5884
$ws.Range("K20").Value = 8078.5884
$ws.Range("M20").Value = -7831.5884
$ws.Range("H132").Value = 90370
$ws.Range("J132").Value = 90370
$ws.Range("L132").Value = 90370
$ws.Range("N132").Value = -100490
$ws.Range("H134").Value = 944
$ws.Range("I134").Value = 944
$ws.Range("K134").Value = 2832
$ws.Range("M134").Value = -297

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 29183.715
$ws.Range("J92").Value = 29183.715
$ws.Range("L92").Value = 29183.715
$ws.Range("N92").Value = -34175.715
$ws.Range("H105").Value = 825.7917
$ws.Range("I105").Value = 844.3043
$ws.Range("J105").Value = 400
$ws.Range("K105").Value = 844.3043
$ws.Range("L105").Value = 400
$ws.Range("M105").Value = 902.6957
$ws.Range("N105").Value = -3894
$ws.Range("H122").Value = 128420.375
$ws.Range("I122").Value = 204534
$ws.Range("J122").Value = 1564.3334
$ws.Range("K122").Value = 613602
$ws.Range("L122").Value = 4693.0002
$ws.Range("M122").Value = -611152
$ws.Range("N122").Value = -9593.0002
$ws.Range("H141").Value = 379356.12
$ws.Range("J141").Value = 379356.12
$ws.Range("L141").Value = 379356.12
$ws.Range("N141").Value = -389716.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3665
$ws.Range("J70").Value = 4499
$ws.Range("L70").Value = 13497
$ws.Range("N70").Value = -14127
$ws.Range("H73").Value = 3665
$ws.Range("J73").Value = 4499
$ws.Range("L73").Value = 13497
$ws.Range("N73").Value = -15681
$ws.Range("H88").Value = 8999.666999999999
$ws.Range("J88").Value = 8999.666999999999
$ws.Range("L88").Value = 26999.001
$ws.Range("N88").Value = -27855.001
$ws.Range("H91").Value = 8999.666999999999
$ws.Range("J91").Value = 8999.666999999999
$ws.Range("L91").Value = 26999.001
$ws.Range("N91").Value = -29963.001
$ws.Range("H113").Value = 1255.1904
$ws.Range("I113").Value = 1066.6666
$ws.Range("J113").Value = 1286.6111
$ws.Range("K113").Value = 3199.9998
$ws.Range("L113").Value = 3859.8333
$ws.Range("M113").Value = -1029.9998
$ws.Range("N113").Value = -8199.8333
$ws.Range("H132").Value = 2520.9333
$ws.Range("I132").Value = 1974.75
$ws.Range("J132").Value = 2719.5454
$ws.Range("K132").Value = 17772.75
$ws.Range("L132").Value = 24475.9086
$ws.Range("M132").Value = -15242.75
$ws.Range("N132").Value = -29535.9086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 400.9091
$ws.Range("I2").Value = 455.3846
$ws.Range("K2").Value = 455.3846
$ws.Range("M2").Value = -342.3846
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 23106.2
$ws.Range("I97").Value = 32418.285
$ws.Range("K97").Value = 32418.285
$ws.Range("M97").Value = -31922.285
$ws.Range("H132").Value = 3972.8572
$ws.Range("I132").Value = 3439.125
$ws.Range("K132").Value = 10317.375
$ws.Range("M132").Value = -7787.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 79163.71000000001
$ws.Range("J16").Value = 839.6
$ws.Range("L16").Value = 839.6
$ws.Range("N16").Value = -1179.6
$ws.Range("H18").Value = 16021.739
$ws.Range("I18").Value = 16021.739
$ws.Range("K18").Value = 16021.739
$ws.Range("M18").Value = -15849.739
$ws.Range("H46").Value = 21161.348
$ws.Range("I46").Value = 48079.89
$ws.Range("K46").Value = 48079.89
$ws.Range("M46").Value = -47891.89
$ws.Range("H132").Value = 3529.4075
$ws.Range("I132").Value = 3011.923
$ws.Range("J132").Value = 4874.8667
$ws.Range("K132").Value = 9035.769
$ws.Range("L132").Value = 14624.6001
$ws.Range("M132").Value = -6505.769
$ws.Range("N132").Value = -19684.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 23810862
$ws.Range("I107").Value = 1777.75
$ws.Range("J107").Value = 55556308
$ws.Range("K107").Value = 5333.25
$ws.Range("L107").Value = 166668924
$ws.Range("M107").Value = -3413.25
$ws.Range("N107").Value = -166672764
$ws.Range("H126").Value = 1575.2858
$ws.Range("I126").Value = 1004.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 3013.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -543.5
$ws.Range("N126").Value = -19940
$ws.Range("H138").Value = 65000
$ws.Range("J138").Value = 65000
$ws.Range("L138").Value = 65000
$ws.Range("N138").Value = -75280
